$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35").Value = "05/01/2026 00:28:34"
$ws.Range("B35").Value = "05/01 00:01"
$ws.Range("C35").Value = "g1 > Economia"
$ws.Range("D35").Value = "Pensando em se demitir? Você não está sozinho: saiba por que tantos brasileiros querem sair do emprego"
$ws.Range("E35").Value = "https://g1.globo.com/trabalho-e-carreira/noticia/2026/01/05/por-que-tantos-brasileiros-querem-sair-do-emprego.ghtml"
$ws.Range("F35").Value = "inflação"
$ws.Range("G35").Value = "dos principais pesos para as empresas, ainda mais com o salário mínimo crescendo acima da &lt;b&gt;inflação&lt;/b&gt;`".`n`"Melhorar benefícios, qualidade de vida e ambiente de trabalho ajuda, mas o que vemos é"
$ws.Range("A36").Value = "05/01/2026 00:28:35"
$ws.Range("B36").Value = "05/01 00:00"
$ws.Range("C36").Value = "g1 > Economia"
$ws.Range("D36").Value = "Saque-aniversário do FGTS: saiba como funciona e veja o calendário para 2026"
$ws.Range("E36").Value = "https://g1.globo.com/economia/noticia/2026/01/05/saque-aniversario-do-fgts-veja-o-calendario-para-2026.ghtml"
$ws.Range("F36").Value = "ldo"
$ws.Range("G36").Value = " Fundo de Garantia por Tempo de Serviço (FGTS) podem retirar, uma vez por ano, parte do sa&lt;b&gt;ldo&lt;/b&gt; das contas ativas e inativas do fundo.`nEssa modalidade é uma alternativa ao modelo tradic"

$ws.Rows.Item(35).AutoFit()
$ws.Rows.Item(36).AutoFit()
